$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl12"
$ws.Range("C2").Value = "Cd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 81.05837566666666
$ws.Range("H2").Value = 243.175127
$ws.Range("I2").Value = 0.3545816884225585
$ws.Range("J2").Value = 0.3545816884225585
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.7633283333333334
$ws.Range("N2").Value = 2.289985
$ws.Range("O2").Value = 0.07247598971834317
$ws.Range("P2").Value = 0.07247598971834318
$ws.Range("Q2").Value = 61.87415480034389
$ws.Range("R2").Value = 556.8673932030949
$ws.Range("S2").Value = 0.02569865880442611
$ws.Range("T2").Value = 0.02569865880442612

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl12"
$ws.Range("C3").Value = "Cd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 81.05837566666666
$ws.Range("H3").Value = 243.175127
$ws.Range("I3").Value = 0.3545816884225585
$ws.Range("J3").Value = 0.3545816884225585
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.451606333333333
$ws.Range("N3").Value = 13.354819
$ws.Range("O3").Value = 0.4226681504613934
$ws.Range("P3").Value = 0.4226681504613934
$ws.Range("Q3").Value = 360.8399784874459
$ws.Range("R3").Value = 3247.559806387012
$ws.Range("S3").Value = 0.1498703864330409
$ws.Range("T3").Value = 0.1498703864330409

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl12"
$ws.Range("C4").Value = "Cd4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 81.05837566666666
$ws.Range("H4").Value = 243.175127
$ws.Range("I4").Value = 0.3545816884225585
$ws.Range("J4").Value = 0.3545816884225585
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.617562
$ws.Range("N4").Value = 13.852686
$ws.Range("O4").Value = 0.438425198465246
$ws.Range("P4").Value = 0.438425198465246
$ws.Range("Q4").Value = 374.2920752601247
$ws.Range("R4").Value = 3368.628677341122
$ws.Range("S4").Value = 0.1554575471188022
$ws.Range("T4").Value = 0.1554575471188022

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cxcl12"
$ws.Range("C5").Value = "Cd4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 81.05837566666666
$ws.Range("H5").Value = 243.175127
$ws.Range("I5").Value = 0.3545816884225585
$ws.Range("J5").Value = 0.3545816884225585
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.699658
$ws.Range("N5").Value = 2.098974
$ws.Range("O5").Value = 0.06643066135501745
$ws.Range("P5").Value = 0.06643066135501746
$ws.Range("Q5").Value = 56.71314100218866
$ws.Range("R5").Value = 510.418269019698
$ws.Range("S5").Value = 0.02355509606628929
$ws.Range("T5").Value = 0.0235550960662893

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl12"
$ws.Range("C6").Value = "Cd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 123.018252
$ws.Range("H6").Value = 369.054756
$ws.Range("I6").Value = 0.5381309351710768
$ws.Range("J6").Value = 0.5381309351710768
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.7633283333333334
$ws.Range("N6").Value = 2.289985
$ws.Range("O6").Value = 0.07247598971834317
$ws.Range("P6").Value = 0.07247598971834318
$ws.Range("Q6").Value = 93.90331726874001
$ws.Range("R6").Value = 845.1298554186601
$ws.Range("S6").Value = 0.03900157212458136
$ws.Range("T6").Value = 0.03900157212458136

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl12"
$ws.Range("C7").Value = "Cd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 123.018252
$ws.Range("H7").Value = 369.054756
$ws.Range("I7").Value = 0.5381309351710768
$ws.Range("J7").Value = 0.5381309351710768
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.451606333333333
$ws.Range("N7").Value = 13.354819
$ws.Range("O7").Value = 0.4226681504613934
$ws.Range("P7").Value = 0.4226681504613934
$ws.Range("Q7").Value = 547.6288297187961
$ws.Range("R7").Value = 4928.659467469163
$ws.Range("S7").Value = 0.227450807074819
$ws.Range("T7").Value = 0.227450807074819

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Cxcl12"
$ws.Range("C8").Value = "Cd4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 123.018252
$ws.Range("H8").Value = 369.054756
$ws.Range("I8").Value = 0.5381309351710768
$ws.Range("J8").Value = 0.5381309351710768
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.617562
$ws.Range("N8").Value = 13.852686
$ws.Range("O8").Value = 0.438425198465246
$ws.Range("P8").Value = 0.438425198465246
$ws.Range("Q8").Value = 568.044405741624
$ws.Range("R8").Value = 5112.399651674617
$ws.Range("S8").Value = 0.2359301620526678
$ws.Range("T8").Value = 0.2359301620526678

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Cxcl12"
$ws.Range("C9").Value = "Cd4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 123.018252
$ws.Range("H9").Value = 369.054756
$ws.Range("I9").Value = 0.5381309351710768
$ws.Range("J9").Value = 0.5381309351710768
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.699658
$ws.Range("N9").Value = 2.098974
$ws.Range("O9").Value = 0.06643066135501745
$ws.Range("P9").Value = 0.06643066135501746
$ws.Range("Q9").Value = 86.070704157816
$ws.Range("R9").Value = 774.636337420344
$ws.Range("S9").Value = 0.03574839391900865
$ws.Range("T9").Value = 0.03574839391900866

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Cxcl12"
$ws.Range("C10").Value = "Cd4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3624666666666667
$ws.Range("H10").Value = 1.0874
$ws.Range("I10").Value = 0.001585573873230423
$ws.Range("J10").Value = 0.001585573873230423
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.7633283333333334
$ws.Range("N10").Value = 2.289985
$ws.Range("O10").Value = 0.07247598971834317
$ws.Range("P10").Value = 0.07247598971834318
$ws.Range("Q10").Value = 0.2766810765555556
$ws.Range("R10").Value = 2.490129689
$ws.Range("S10").Value = 0.0001149160357339217
$ws.Range("T10").Value = 0.0001149160357339217

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Cxcl12"
$ws.Range("C11").Value = "Cd4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3624666666666667
$ws.Range("H11").Value = 1.0874
$ws.Range("I11").Value = 0.001585573873230423
$ws.Range("J11").Value = 0.001585573873230423
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.451606333333333
$ws.Range("N11").Value = 13.354819
$ws.Range("O11").Value = 0.4226681504613934
$ws.Range("P11").Value = 0.4226681504613934
$ws.Range("Q11").Value = 1.613558908955556
$ws.Range("R11").Value = 14.5220301806
$ws.Range("S11").Value = 0.0006701715764182109
$ws.Range("T11").Value = 0.0006701715764182109

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Cxcl12"
$ws.Range("C12").Value = "Cd4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3624666666666667
$ws.Range("H12").Value = 1.0874
$ws.Range("I12").Value = 0.001585573873230423
$ws.Range("J12").Value = 0.001585573873230423
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.617562
$ws.Range("N12").Value = 13.852686
$ws.Range("O12").Value = 0.438425198465246
$ws.Range("P12").Value = 0.438425198465246
$ws.Range("Q12").Value = 1.673712306266667
$ws.Range("R12").Value = 15.0634107564
$ws.Range("S12").Value = 0.0006951555400523573
$ws.Range("T12").Value = 0.0006951555400523573

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Cxcl12"
$ws.Range("C13").Value = "Cd4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3624666666666667
$ws.Range("H13").Value = 1.0874
$ws.Range("I13").Value = 0.001585573873230423
$ws.Range("J13").Value = 0.001585573873230423
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.699658
$ws.Range("N13").Value = 2.098974
$ws.Range("O13").Value = 0.06643066135501745
$ws.Range("P13").Value = 0.06643066135501746
$ws.Range("Q13").Value = 0.2536027030666667
$ws.Range("R13").Value = 2.2824243276
$ws.Range("S13").Value = 0.0001053307210259336
$ws.Range("T13").Value = 0.0001053307210259336

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Cxcl12"
$ws.Range("C14").Value = "Cd4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 24.16373066666667
$ws.Range("H14").Value = 72.491192
$ws.Range("I14").Value = 0.1057018025331343
$ws.Range("J14").Value = 0.1057018025331344
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.7633283333333334
$ws.Range("N14").Value = 2.289985
$ws.Range("O14").Value = 0.07247598971834317
$ws.Range("P14").Value = 0.07247598971834318
$ws.Range("Q14").Value = 18.44486025690222
$ws.Range("R14").Value = 166.00374231212
$ws.Range("S14").Value = 0.007660842753601784
$ws.Range("T14").Value = 0.007660842753601787

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Cxcl12"
$ws.Range("C15").Value = "Cd4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 24.16373066666667
$ws.Range("H15").Value = 72.491192
$ws.Range("I15").Value = 0.1057018025331343
$ws.Range("J15").Value = 0.1057018025331344
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.451606333333333
$ws.Range("N15").Value = 13.354819
$ws.Range("O15").Value = 0.4226681504613934
$ws.Range("P15").Value = 0.4226681504613934
$ws.Range("Q15").Value = 107.5674164726942
$ws.Range("R15").Value = 968.1067482542479
$ws.Range("S15").Value = 0.04467678537711531
$ws.Range("T15").Value = 0.04467678537711532

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Cxcl12"
$ws.Range("C16").Value = "Cd4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 24.16373066666667
$ws.Range("H16").Value = 72.491192
$ws.Range("I16").Value = 0.1057018025331343
$ws.Range("J16").Value = 0.1057018025331344
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 4.617562
$ws.Range("N16").Value = 13.852686
$ws.Range("O16").Value = 0.438425198465246
$ws.Range("P16").Value = 0.438425198465246
$ws.Range("Q16").Value = 111.5775245046347
$ws.Range("R16").Value = 1004.197720541712
$ws.Range("S16").Value = 0.04634233375372367
$ws.Range("T16").Value = 0.04634233375372367

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Cxcl12"
$ws.Range("C17").Value = "Cd4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 24.16373066666667
$ws.Range("H17").Value = 72.491192
$ws.Range("I17").Value = 0.1057018025331343
$ws.Range("J17").Value = 0.1057018025331344
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.699658
$ws.Range("N17").Value = 2.098974
$ws.Range("O17").Value = 0.06643066135501745
$ws.Range("P17").Value = 0.06643066135501746
$ws.Range("Q17").Value = 16.90634747077867
$ws.Range("R17").Value = 152.157127237008
$ws.Range("S17").Value = 0.007021840648693572
$ws.Range("T17").Value = 0.007021840648693575
